# Updated cryptos list on Tue May  9 09:16:13 UTC 2023 with GitHub Actions
# Refresh Price (D) and Volume(1h) (E) columns for the cryptos table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.622.92"
$ws.Range("E2").Value = "  -0.95%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.844.46"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4221"
$ws.Range("E7").Value = "  -2.90%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3636"
$ws.Range("E8").Value = "  -1.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.26"
$ws.Range("E9").Value = "  +1.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07262"
$ws.Range("E10").Value = "  -2.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8888"
$ws.Range("E11").Value = "  -5.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.59"
$ws.Range("E12").Value = "  -3.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.862.27"
$ws.Range("E13").Value = "  +0.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.563"
$ws.Range("E14").Value = "  -2.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.332"
$ws.Range("E15").Value = "  -1.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06862"
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.004"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "78.82"
$ws.Range("E18").Value = "  -3.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008823"
$ws.Range("E19").Value = "  -2.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.43"
$ws.Range("E21").Value = "  -2.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.606.06"
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.977"
$ws.Range("E23").Value = "  -2.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.53"
$ws.Range("E24").Value = "  -4.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.065.45"
$ws.Range("E25").Value = "  -1.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.923"
$ws.Range("E26").Value = "  -4.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.84"
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.48"
$ws.Range("E28").Value = "  +0.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "119.63"
$ws.Range("E29").Value = "  +5.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.237"
$ws.Range("E30").Value = "  -2.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.847"
$ws.Range("E31").Value = "  +6.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08893"
$ws.Range("E32").Value = "  -0.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7776"
$ws.Range("E33").Value = "  -2.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.564"
$ws.Range("E34").Value = "  -5.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.946"
$ws.Range("E35").Value = "  -1.99%  "
$ws.Range("E36").Value = "  -6.24%  "
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05398"
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.098"
$ws.Range("E39").Value = "  -1.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01926"
$ws.Range("E40").Value = "  -1.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.804"
$ws.Range("E41").Value = "  -3.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.853"
$ws.Range("E42").Value = "  -2.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5061"
$ws.Range("E43").Value = "  -3.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1648"
$ws.Range("E44").Value = "  -1.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.263"
$ws.Range("E45").Value = "  -5.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06599"
$ws.Range("E46").Value = "  -1.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.30"
$ws.Range("E47").Value = "  -2.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4697"
$ws.Range("E48").Value = "  -3.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.47"
$ws.Range("E49").Value = "  -2.22%  "
$ws.Range("E50").Value = "  -0.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.630"
$ws.Range("E51").Value = "  -2.71%  "
